$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.314.94"
$ws.Range("E2").Value = "  +6.59%  "
$ws.Range("D3").Value = "3.000.33"
$ws.Range("E3").Value = "  +4.14%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'583.97"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").Value = "'154.49"
$ws.Range("E6").Value = "  +7.76%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.997.12"
$ws.Range("E8").Value = "  +3.96%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").Value = "'7.03"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  +4.84%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "'34.08"
$ws.Range("E14").Value = "  +6.88%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "65.234.93"
$ws.Range("E16").Value = "  +6.35%  "
$ws.Range("D17").Value = "3.492.53"
$ws.Range("E17").Value = "  +4.09%  "
$ws.Range("D18").Value = "'6.95"
$ws.Range("E18").Value = "  +5.92%  "
$ws.Range("D19").Value = "2.997.31"
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("D20").Value = "'451.08"
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("D21").Value = "'13.73"
$ws.Range("E21").Value = "  +5.35%  "
$ws.Range("D22").Value = "'0.682"
$ws.Range("E22").Value = "  +4.46%  "
$ws.Range("D23").Value = "'7.34"
$ws.Range("E23").Value = "  +7.64%  "
$ws.Range("D24").Value = "'81.25"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "'12.49"
$ws.Range("E25").Value = "  +5.79%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  +11.57%  "
$ws.Range("D27").Value = "'10.77"
$ws.Range("E27").Value = "  +7.87%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'2.43"
$ws.Range("E29").Value = "  +18.49%  "
$ws.Range("D30").Value = "'7.82"
$ws.Range("E30").Value = "  +12.06%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0000104"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.60"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("D34").Value = "'26.92"
$ws.Range("E34").Value = "  +5.78%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "'0.983"
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").Value = "'5.80"
$ws.Range("E37").Value = "  +7.92%  "
$ws.Range("D38").Value = "'2.11"
$ws.Range("E38").Value = "  +8.75%  "
$ws.Range("D39").Value = "'46.27"
$ws.Range("E39").Value = "  +18.46%  "
$ws.Range("D40").Value = "'49.12"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "'2.94"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").Value = "'0.303"
$ws.Range("E42").Value = "  +13.61%  "
$ws.Range("D43").Value = "'0.121"
$ws.Range("E43").Value = "  +6.55%  "
$ws.Range("D44").Value = "'8.41"
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("D45").Value = "'385.81"
$ws.Range("E45").Value = "  +12.92%  "
$ws.Range("D46").Value = "2.766.86"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").Value = "'0.0350"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("D48").Value = "'134.75"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'23.19"
$ws.Range("E50").Value = "  +7.75%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +2.78%  "
